{"js": "// Guidance Admission Slip: fill in the templated blanks (student name,\n// grade/section, teacher name, visit date/time range, counselor & teacher\n// signature lines). The document repeats the same letter twice (a\n// \"GUIDANCE COPY\" and a \"TEACHER'S COPY\" further down), so the same\n// placeholder text occurs twice; replacements are scoped per-paragraph so\n// each occurrence is only touched once, in the right place.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// [paragraphIndex, searchText, replacementText]\nconst edits = [\n  [5, \"John Vincent\", \"acsasca\"],\n  [5, \"12-ambot\", \"scacscasc\"],\n  [7, \"example teacher\", \"ascascas\"],\n  [8, \"  09:01:00   \", \"  17:19:00   \"],\n  [8, \"19:57:00\", \"18:20:00\"],\n  [12, \"Example Counselor________________\", \"ascascascas________________\"],\n  [15, \"example teacher_________\", \"ascascas_________\"],\n  [26, \"John Vincent\", \"acsasca\"],\n  [26, \"12-ambot\", \"scacscasc\"],\n  [28, \"example teacher\", \"ascascas\"],\n  [29, \"  09:01:00   \", \"  17:19:00   \"],\n  [29, \"19:57:00\", \"18:20:00\"],\n  [33, \"Example Counselor________________\", \"ascascascas________________\"],\n  [36, \"example teacher_________\", \"ascascas_________\"],\n];\n\nfor (const [paraIndex, findText, replaceText] of edits) {\n  const paragraph = paragraphs.items[paraIndex];\n  const found = paragraph.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  found.load(\"items\");\n  await context.sync();\n\n  for (const range of found.items) {\n    range.insertText(replaceText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Guidance Admission Slip: fill in the templated blanks (student name,\n# grade/section, teacher name, visit date/time range, counselor & teacher\n# signature lines). The document repeats the same letter twice (a\n# \"GUIDANCE COPY\" and a \"TEACHER'S COPY\" further down), so the same\n# placeholder text occurs twice; replacements are scoped to a specific\n# paragraph (via Paragraphs.Item, 1-based) so each occurrence is only\n# touched once, in the right place.\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n\n$edits = @(\n    @{ Para = 6;  Find = \"John Vincent\"; Replace = \"acsasca\" },\n    @{ Para = 6;  Find = \"12-ambot\"; Replace = \"scacscasc\" },\n    @{ Para = 8;  Find = \"example teacher\"; Replace = \"ascascas\" },\n    @{ Para = 9;  Find = \"  09:01:00   \"; Replace = \"  17:19:00   \" },\n    @{ Para = 9;  Find = \"19:57:00\"; Replace = \"18:20:00\" },\n    @{ Para = 13; Find = \"Example Counselor________________\"; Replace = \"ascascascas________________\" },\n    @{ Para = 16; Find = \"example teacher_________\"; Replace = \"ascascas_________\" },\n    @{ Para = 27; Find = \"John Vincent\"; Replace = \"acsasca\" },\n    @{ Para = 27; Find = \"12-ambot\"; Replace = \"scacscasc\" },\n    @{ Para = 29; Find = \"example teacher\"; Replace = \"ascascas\" },\n    @{ Para = 30; Find = \"  09:01:00   \"; Replace = \"  17:19:00   \" },\n    @{ Para = 30; Find = \"19:57:00\"; Replace = \"18:20:00\" },\n    @{ Para = 34; Find = \"Example Counselor________________\"; Replace = \"ascascascas________________\" },\n    @{ Para = 37; Find = \"example teacher_________\"; Replace = \"ascascas_________\" }\n)\n\nforeach ($edit in $edits) {\n    $range = $d.Paragraphs.Item($edit.Para).Range\n    $find = $range.Find\n    $find.Text = $edit.Find\n    $find.Replacement.Text = $edit.Replace\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Wrap = 0\n    $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, $wdReplaceOne)\n}\n"}
